{"js": "// Replace the division-problem answers in the table cells.\n// Each entry is [oldText, newText]; every oldText occurs exactly once\n// in the document (inside a table cell paragraph), so a plain\n// case-sensitive whole-match search is unambiguous.\nconst replacements = [\n  [\"63\u00f77=9, 0\", \"39\u00f72=19, 1\"],\n  [\"94\u00f72=47, 0\", \"97\u00f79=10, 7\"],\n  [\"87\u00f72=43, 1\", \"17\u00f74=4, 1\"],\n  [\"52\u00f78=6, 4\", \"97\u00f77=13, 6\"],\n  [\"31\u00f77=4, 3\", \"40\u00f72=20, 0\"],\n  [\"80\u00f73=26, 2\", \"31\u00f79=3, 4\"],\n  [\"66\u00f72=33, 0\", \"56\u00f78=7, 0\"],\n  [\"74\u00f75=14, 4\", \"84\u00f72=42, 0\"],\n  [\"58\u00f75=11, 3\", \"21\u00f74=5, 1\"],\n  [\"43\u00f75=8, 3\", \"34\u00f79=3, 7\"],\n  [\"62\u00f79=6, 8\", \"44\u00f78=5, 4\"],\n  [\"45\u00f74=11, 1\", \"29\u00f72=14, 1\"],\n  [\"34\u00f77=4, 6\", \"71\u00f73=23, 2\"],\n  [\"93\u00f75=18, 3\", \"16\u00f77=2, 2\"],\n  [\"73\u00f79=8, 1\", \"94\u00f78=11, 6\"],\n  [\"47\u00f74=11, 3\", \"29\u00f73=9, 2\"],\n  [\"18\u00f79=2, 0\", \"39\u00f77=5, 4\"],\n  [\"75\u00f76=12, 3\", \"55\u00f78=6, 7\"],\n  [\"58\u00f78=7, 2\", \"11\u00f76=1, 5\"],\n  [\"35\u00f73=11, 2\", \"35\u00f76=5, 5\"],\n  [\"96\u00f79=10, 6\", \"23\u00f79=2, 5\"],\n  [\"86\u00f75=17, 1\", \"28\u00f77=4, 0\"],\n  [\"46\u00f78=5, 6\", \"13\u00f79=1, 4\"],\n  [\"92\u00f75=18, 2\", \"57\u00f77=8, 1\"],\n  [\"89\u00f72=44, 1\", \"15\u00f73=5, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem answers in the table cells.\n# Each pair is the old (Find What) text and the new (Replace With)\n# text; every old text occurs exactly once in the document, so a\n# simple Find/Replace (MatchCase, whole match) is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"63\u00f77=9, 0\", \"39\u00f72=19, 1\"),\n    @(\"94\u00f72=47, 0\", \"97\u00f79=10, 7\"),\n    @(\"87\u00f72=43, 1\", \"17\u00f74=4, 1\"),\n    @(\"52\u00f78=6, 4\", \"97\u00f77=13, 6\"),\n    @(\"31\u00f77=4, 3\", \"40\u00f72=20, 0\"),\n    @(\"80\u00f73=26, 2\", \"31\u00f79=3, 4\"),\n    @(\"66\u00f72=33, 0\", \"56\u00f78=7, 0\"),\n    @(\"74\u00f75=14, 4\", \"84\u00f72=42, 0\"),\n    @(\"58\u00f75=11, 3\", \"21\u00f74=5, 1\"),\n    @(\"43\u00f75=8, 3\", \"34\u00f79=3, 7\"),\n    @(\"62\u00f79=6, 8\", \"44\u00f78=5, 4\"),\n    @(\"45\u00f74=11, 1\", \"29\u00f72=14, 1\"),\n    @(\"34\u00f77=4, 6\", \"71\u00f73=23, 2\"),\n    @(\"93\u00f75=18, 3\", \"16\u00f77=2, 2\"),\n    @(\"73\u00f79=8, 1\", \"94\u00f78=11, 6\"),\n    @(\"47\u00f74=11, 3\", \"29\u00f73=9, 2\"),\n    @(\"18\u00f79=2, 0\", \"39\u00f77=5, 4\"),\n    @(\"75\u00f76=12, 3\", \"55\u00f78=6, 7\"),\n    @(\"58\u00f78=7, 2\", \"11\u00f76=1, 5\"),\n    @(\"35\u00f73=11, 2\", \"35\u00f76=5, 5\"),\n    @(\"96\u00f79=10, 6\", \"23\u00f79=2, 5\"),\n    @(\"86\u00f75=17, 1\", \"28\u00f77=4, 0\"),\n    @(\"46\u00f78=5, 6\", \"13\u00f79=1, 4\"),\n    @(\"92\u00f75=18, 2\", \"57\u00f77=8, 1\"),\n    @(\"89\u00f72=44, 1\", \"15\u00f73=5, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $null, $find.Replacement.Text, 2)\n}\n\n$d.Save()\n"}
